$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: becomes the start of a new block (style flips to the s=6/7 group) ---
# Copy the cell formatting (not values) from row 26 (an existing "start" row in the
# s=6/7 style group) onto row 28 so B28/C28/D28/E28 keep their current values but
# pick up the new style, and A28 becomes an empty styled cell ready for its value.
$ws.Range("A26:E26").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Rows.Item(28).RowHeight = 43.2

# --- Rows 29-31: brand new rows ---
# Row 29 uses the s=4/5 "start" style (like row 25, which has a value in column A).
$ws.Range("A25:E25").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Rows.Item(29).RowHeight = 43.2

# Rows 30-31 use the s=4/5 "continuation" style (like row 9, which has no value in
# column A), so we only paste formats onto columns B:E and leave A untouched (no
# A cell should exist on these rows).
$ws.Range("B9:E9").Copy()
$ws.Range("B30:E30").PasteSpecial(-4122)
$ws.Rows.Item(30).RowHeight = 31.8

$ws.Range("B9:E9").Copy()
$ws.Range("B31:E31").PasteSpecial(-4122)

# --- Values, set in the same order the original file's translators used so new
# shared strings land at the expected indices (96-106): filename (EN), the 3
# English lines, filename (RU), the 3 Russian lines, then the 3 "encoded" lines ---
$ws.Range("A28").Value = "SCRIPT/G01P04A/um2502.ssb"

$ws.Range("C29").Value = " Excuse me?![K] The grand master of\nall things bad?!"
$ws.Range("C30").Value = " I...I really can\'t tell you...[K]\nPlease, I just can\'t…"
$ws.Range("C31").Value = " I\'m so sorry…"

$ws.Range("A29").Value = "SCRIPT/G01P04A/us0103.ssb"

$ws.Range("D29").Value = " Что, простите?![K] Гранд мастер\nвсего самого плохого?!"
$ws.Range("D30").Value = " Я... Я правда не могу вам\nсказать...[K] Прошу, я совсем не могу..."
$ws.Range("D31").Value = " Мне очень жаль..."

$ws.Range("E29").Value = " Œóï, ðñïòóéóå?![K] Ãñàîä íàòóåñ\nâòåãï òàíïãï ðìïöïãï?!"
$ws.Range("E30").Value = " Ÿ... Ÿ ðñàâäà îå íïãô âàí\nòëàèàóû...[K] Ðñïšô, ÿ òïâòåí îå íïãô..."
$ws.Range("E31").Value = " Íîå ïœåîû çàìû..."

# --- B column numeric values (unchanged numbering style) ---
$ws.Range("B29").Value = 154
$ws.Range("B30").Value = 157
$ws.Range("B31").Value = 160

# --- View state: move the selection to match the post-edit workbook ---
$ws.Range("D9").Select()
